$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per diff (D/E columns forced to Text format so
# numeric-looking strings such as prices and percentages are preserved exactly,
# matching the original inlineStr text cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.17'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '36.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-4.08%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.112'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.56%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07706'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.64%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.388'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.72%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.298'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.49%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.842'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.41%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.945'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.32%'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9213'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.03%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1115'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-8.25%'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1843'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.92%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08796'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.85%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03340'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.74%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09518'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.95%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001388'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.64%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006179'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '8.12%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.362'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-4.49%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3436'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.52%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.338'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '20.25%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1291'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.97%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-10.80%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04339'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.54%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.91%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.26%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001329'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '8.92%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02075'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.95%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04938'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.51%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007506'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.61%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1348'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.97%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008401'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.79%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002068'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.62%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008406'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.63%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006295'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.84%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002852'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-14.86%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '40.95%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.12%'
